# "Fruta / hortaliza, semanal"
#
# A new weekly price-report row is inserted at row 362 of the sheet
# (pushing the existing rows 362-379 down to 363-380), and populated with
# the new observation. The worksheet dimension / row count grows from
# A1:R379 to A1:R380 automatically as a result of the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 362, shifting rows 362:379
# down to 363:380.
$ws.Rows.Item(362).Insert()

# Populate the newly inserted row 362 with the new weekly observation.
$ws.Cells.Item(362, 1).Value = 5
$ws.Cells.Item(362, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(362, 3).Value = "Maule"
$ws.Cells.Item(362, 4).Value = 44753
$ws.Cells.Item(362, 5).Value = 7
$ws.Cells.Item(362, 6).Value = 100112032
$ws.Cells.Item(362, 7).Value = "Zapallo italiano"
$ws.Cells.Item(362, 8).Value = "Sin especificar"
$ws.Cells.Item(362, 9).Value = "Primera"
$ws.Cells.Item(362, 10).Value = 300
$ws.Cells.Item(362, 11).Value = 11000
$ws.Cells.Item(362, 12).Value = 11000
$ws.Cells.Item(362, 13).Value = 11000
$ws.Cells.Item(362, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(362, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(362, 16).Value = 220
$ws.Cells.Item(362, 17).Value = 50
$ws.Cells.Item(362, 18).Value = "Hortaliza"
